# Apply crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.305.20'
$ws.Range("E2").Value = '  -4.19%  '
$ws.Range("D3").Value = '3.674.25'
$ws.Range("E3").Value = '  -5.18%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.44'
$ws.Range("E5").Value = '  -2.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.11'
$ws.Range("E6").Value = '  +4.94%  '
$ws.Range("D7").Value = '3.671.30'
$ws.Range("E7").Value = '  -5.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.630'
$ws.Range("E8").Value = '  -6.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.713'
$ws.Range("E10").Value = '  -4.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.161'
$ws.Range("E11").Value = '  -9.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.36'
$ws.Range("E12").Value = '  +2.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000294'
$ws.Range("E13").Value = '  -9.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.64'
$ws.Range("E14").Value = '  -7.94%  '
$ws.Range("D15").Value = '4.252.57'
$ws.Range("E15").Value = '  -5.16%  '
$ws.Range("D16").Value = '3.668.44'
$ws.Range("E16").Value = '  -5.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.34'
$ws.Range("E17").Value = '  -8.99%  '
$ws.Range("E18").Value = '  -2.18%  '
$ws.Range("E19").Value = '  -7.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.78'
$ws.Range("E20").Value = '  -8.51%  '
$ws.Range("D21").Value = '68.166.03'
$ws.Range("E21").Value = '  -4.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '409.26'
$ws.Range("E22").Value = '  -7.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.57'
$ws.Range("E23").Value = '  -5.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.56'
$ws.Range("E24").Value = '  -6.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.00'
$ws.Range("E25").Value = '  -9.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.65'
$ws.Range("E26").Value = '  -9.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.79'
$ws.Range("E27").Value = '  -7.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.88'
$ws.Range("E28").Value = '  -3.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.56'
$ws.Range("E30").Value = '  -9.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.75'
$ws.Range("E31").Value = '  -7.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.30'
$ws.Range("E32").Value = '  -16.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.38'
$ws.Range("E33").Value = '  -8.87%  '
$ws.Range("E34").Value = '  -6.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '64.58'
$ws.Range("E35").Value = '  -6.53%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '43.18'
$ws.Range("E36").Value = '  -10.80%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '601.06'
$ws.Range("E37").Value = '  -6.13%  '
$ws.Range("D38").Value = '0.0₃0885'
$ws.Range("E38").Value = '  -10.74%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.399'
$ws.Range("E40").Value = '  -9.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.137'
$ws.Range("E42").Value = '  -6.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.02'
$ws.Range("E43").Value = '  -7.83%  '
$ws.Range("E44").Value = '  -6.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0441'
$ws.Range("E45").Value = '  -6.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.81'
$ws.Range("E46").Value = '  -10.35%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.135'
$ws.Range("E47").Value = '  -6.73%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.72'
$ws.Range("E48").Value = '  -6.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.96'
$ws.Range("E49").Value = '  -12.64%  '
$ws.Range("D50").Value = '2.715.83'
$ws.Range("E50").Value = '  -8.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.10'
$ws.Range("E51").Value = '  -6.90%  '
